$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.112.82'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.31%  '

$ws.Range("D3").Value = '''1.834.65'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.29%  '

$ws.Range("D4").Value = '''1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.40%  '

$ws.Range("D5").Value = '''242.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.37%  '

$ws.Range("D6").Value = '''0.6172'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.64%  '

$ws.Range("D7").Value = '''1.003'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.22%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = '''0.2965'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.41%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '''0.07397'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.69%  '

$ws.Range("D10").Value = '''22.97'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.05%  '

$ws.Range("D11").Value = '''0.07676'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.30%  '

$ws.Range("D12").Value = '''1.824.06'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.29%  '

$ws.Range("D13").Value = '''5.000'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.07%  '

$ws.Range("D14").Value = '''0.6706'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.78%  '

$ws.Range("D15").Value = '''82.77'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.38%  '

$ws.Range("D16").Value = '''0.000009024'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.80%  '

$ws.Range("D17").Value = '''5.885'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.46%  '

$ws.Range("D18").Value = '''29.095.55'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.28%  '

$ws.Range("D19").Value = '''2.080.18'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.18%  '

$ws.Range("D20").Value = '''236.16'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.50%  '

$ws.Range("D21").Value = '''12.65'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.73%  '

$ws.Range("E22").Value = '  +0.41%  '

$ws.Range("D23").Value = '''7.178'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.12%  '

$ws.Range("D24").Value = '''1.007'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.50%  '

$ws.Range("D25").Value = '''159.26'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.19%  '

$ws.Range("D26").Value = '''0.1431'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.09%  '

$ws.Range("D27").Value = '''8.495'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.13%  '

$ws.Range("D28").Value = '''17.82'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.34%  '

$ws.Range("D29").Value = '''1.499'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.45%  '

$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").Value = '''4.136'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.51%  '

$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").Value = '''0.05563'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.93%  '

$ws.Range("D32").Value = '''4.101'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.25%  '

$ws.Range("D33").Value = '''1.213'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.54%  '

$ws.Range("D34").Value = '''0.7468'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.25%  '

$ws.Range("D35").Value = '''1.847'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.02%  '

$ws.Range("D36").Value = '''1.138'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.35%  '

$ws.Range("D37").Value = '''2.653'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.38%  '

$ws.Range("D38").Value = '''2.791'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.60%  '

$ws.Range("D39").Value = '''0.01779'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.12%  '

$ws.Range("D40").Value = '''1.207.59'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.47%  '

$ws.Range("D41").Value = '''6.367'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.06%  '

$ws.Range("D42").Value = '''0.9025'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.27%  '

$ws.Range("D43").Value = '''1.002'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.08%  '

$ws.Range("D44").Value = '''101.09'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.42%  '

$ws.Range("D45").Value = '''1.974.03'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.31%  '

$ws.Range("E46").Value = '  -0.27%  '

$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '''0.5108'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.33%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '''0.00000000121'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.45%  '

$ws.Range("D49").Value = '''0.4048'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.44%  '

$ws.Range("D50").Value = '''9.141'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.66%  '

$ws.Range("D51").Value = '''0.05825'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.45%  '
